$wb = $excel.ActiveWorkbook
$originalActiveSheet = $wb.ActiveSheet

# ---------------------------------------------------------------------
# Service Contacts sheet: widen column A and move the selection to D3
# ---------------------------------------------------------------------
$wsSC = $wb.Worksheets.Item("Service Contacts")
$wsSC.Columns.Item(1).ColumnWidth = 13.666666666666666
[void]$wsSC.Activate()
[void]$wsSC.Range("D3").Select()

# ---------------------------------------------------------------------
# Practitioners sheet: add a new practitioner row, widen a few columns
# and move the selection to column G
# ---------------------------------------------------------------------
$wsP = $wb.Worksheets.Item("Practitioners")

$wsP.Range("A6").Value = "PHN999:NFP02"
$wsP.Range("B6").Value = "P01"
$wsP.Range("C6").Value = 8
$wsP.Range("D6").Value = 1
$wsP.Range("E6").Value = 1973
$wsP.Range("F6").Value = 2
$wsP.Range("G6").Value = 1
$wsP.Range("H6").Value = 1
$wsP.Range("I6").Value = "tag1"

$wsP.Columns.Item(1).ColumnWidth = 13.833333333333332
$wsP.Columns.Item(3).ColumnWidth = 12.166666666666666
$wsP.Columns.Item(6).ColumnWidth = 12.0

[void]$wsP.Activate()
[void]$wsP.Range("G:G").Select()

# Restore original active sheet/tab so the workbook-level active tab
# is unchanged, matching the source edit.
[void]$originalActiveSheet.Activate()
